$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.250.66'
$ws.Range("E2").Value = '  +2.71%  '
$ws.Range("D3").Value = '1.718.35'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.55'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -1.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2620'
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06194'
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").Value = '1.717.44'
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07074'
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.41'
$ws.Range("E12").Value = '  +3.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5954'
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.422'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.15'
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '26.261.38'
$ws.Range("E18").Value = '  +2.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006798'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '1.937.56'
$ws.Range("E21").Value = '  +3.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.536'
$ws.Range("E22").Value = '  +2.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.725'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.271'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.84'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.16'
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.401'
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.760'
$ws.Range("E28").Value = '  +2.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.95'
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.967'
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.680'
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07753'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04470'
$ws.Range("E33").Value = '  +5.45%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.616'
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9749'
$ws.Range("E35").Value = '  +2.42%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6172'
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9240'
$ws.Range("E37").Value = '  +7.21%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '113.75'
$ws.Range("E38").Value = '  +17.18%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.434'
$ws.Range("E39").Value = '  -6.41%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.920'
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01481'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.597'
$ws.Range("E43").Value = '  +15.43%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3821'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1177'
$ws.Range("E45").Value = '  +4.65%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.286'
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05269'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.37'
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.760'
$ws.Range("E49").Value = '  +5.83%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3373'
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.217'
$ws.Range("E51").Value = '  +1.34%  '
